# Append the "Iteration 2:" section after the existing "Collision Detection"
# paragraph (end of the "Iteration 1:" list):
#   - a blank paragraph
#   - a bold/underlined "Iteration 2:" heading paragraph
#   - three plain list paragraphs: "1. Enemies spawning", "2. AI",
#     "3. Item spawning"

$d = $word.ActiveDocument

# Locate the "Collision Detection" run so we can anchor the insertion right
# after the paragraph that contains it (the last paragraph of "Iteration 1:").
$findRange = $d.Content
$found = $findRange.Find.Execute("Collision Detection", $true, $false, $false,
                                  $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Collision Detection' text"
}

$anchorPara = $findRange.Paragraphs.Last
$insertionPoint = $d.Range($anchorPara.Range.End, $anchorPara.Range.End)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParagraphsXml = (
    "<w:p $wNs/>" +
    "<w:p $wNs><w:r><w:rPr><w:b/><w:u w:val=`"single`"/></w:rPr><w:t>Iteration 2:</w:t></w:r></w:p>" +
    "<w:p $wNs><w:r><w:t>1. Enemies spawning</w:t></w:r></w:p>" +
    "<w:p $wNs><w:r><w:t>2. AI</w:t></w:r></w:p>" +
    "<w:p $wNs><w:r><w:t>3. Item spawning</w:t></w:r></w:p>"
)

$insertionPoint.InsertXML($newParagraphsXml)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
